# Bug fixes and more error messages
# Applies the "Ongoing Games" sheet updates described by the commit:
#  - game on row 3 (Omaha @ James Madison) advances: score, field position,
#    next play type/number and the "Number Submitted" flag
#  - game on row 4 gets its home user's (misspelled) tag corrected in both
#    the "Home User" and "Waiting On" columns, the defensive number and the
#    "Number Submitted" flag are updated
#  - the trailing block of empty placeholder rows (9:11) is trimmed
#  - selection / scroll position left where the editor's cursor ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ongoing Games")

# --- Row 3: Omaha vs James Madison -----------------------------------
$ws.Range("S3").Value = 28
$ws.Range("AD3").Value = "James Madison 35"
$ws.Range("AI3").Value = "KICKOFF"
$ws.Range("AJ3").Value = 665
$ws.Range("AP3").Value = "NO"

# --- Row 4: corrected Discord tag --------------------------------------
$ws.Range("D4").Value = "Buttersqauch#3186"
$ws.Range("AH4").Value = "Buttersqauch#3186"
$ws.Range("AK4").Value = 55
$ws.Range("AP4").Value = "YES"

# --- Trim the trailing empty rows (9:11) -------------------------------
$ws.Rows("9:11").Delete()

# --- Leave the view/selection where the edits ended -------------------
$null = $ws.Range("D4").Select()
